# template_file_format.xlsx update
# - re-sorts/rebuilds the data rows (code, descr, lang_code) grouped by
#   code (html, json, txt, xml) then language, and adds a new "con"
#   language row for the html code (with a "html File" label typo kept
#   as authored).
# - turns on AutoFilter for the header row.
# - re-selects cell K10 to match the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target rows (code, descr, lang_code) in final on-disk order.
$data = @(
  @("html", "html file", "ara"),
  @("html", "html File", "con"),
  @("html", "html file", "eng"),
  @("html", "html file", "fra"),
  @("html", "html file", "hin"),
  @("html", "html file", "kan"),
  @("html", "html file", "tam"),
  @("json", "Json File", "ara"),
  @("json", "Json File", "eng"),
  @("json", "Json File", "fra"),
  @("json", "Json File", "hin"),
  @("json", "Json File", "kan"),
  @("json", "Json File", "tam"),
  @("txt", "Text File", "ara"),
  @("txt", "Text File", "eng"),
  @("txt", "Text File", "fra"),
  @("txt", "Text File", "hin"),
  @("txt", "Text File", "kan"),
  @("txt", "Text File", "tam"),
  @("xml", "XML File", "ara"),
  @("xml", "XML File", "eng"),
  @("xml", "XML File", "fra"),
  @("xml", "XML File", "hin"),
  @("xml", "XML File", "kan"),
  @("xml", "XML File", "tam")
)

$crDate = 45337.275843611111
$dateFmt = $ws.Cells.Item(2, 6).NumberFormat

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  if ($r -eq 3) {
    # the new "con" row: the authored file's shared-string table shows
    # lang_code (C) allocated before descr (B), so poke C first here too.
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 2).Value = $row[1]
  } else {
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
  }
  $ws.Cells.Item($r, 4).Value = "t"
  $ws.Cells.Item($r, 5).Value = "admin"
  $ws.Cells.Item($r, 6).Value = $crDate
  $ws.Cells.Item($r, 6).NumberFormat = $dateFmt
  $ws.Cells.Item($r, 9).Value = "f"
  $r = $r + 1
}

# Row 3 (new "con" row) picked up an extra number-format stamp on the
# cr_by cell in the authored file; replicate it.
$ws.Cells.Item(3, 5).NumberFormat = $dateFmt

# Turn on the header AutoFilter and keep the hidden _FilterDatabase name
# that Excel writes alongside it.
$ws.Range("A1:J1").AutoFilter()
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=template_file_format!`$A`$1:`$J`$1")
$fdb.Visible = $false

# Restore the saved selection.
$ws.Range("K10").Select()
